$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp7"
$ws.Cells.Item(2, 3).Value = "Acvr1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1301303333333333
$ws.Cells.Item(2, 8).Value = 0.390391
$ws.Cells.Item(2, 9).Value = 0.2986126887311651
$ws.Cells.Item(2, 10).Value = 0.2986126887311651
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.489835333333334
$ws.Cells.Item(2, 14).Value = 25.469506
$ws.Cells.Item(2, 15).Value = 0.2075776945087381
$ws.Cells.Item(2, 16).Value = 0.2075776945087381
$ws.Cells.Item(2, 17).Value = 1.104785101871778
$ws.Cells.Item(2, 18).Value = 9.943065916845999
$ws.Cells.Item(2, 19).Value = 0.0619853334778707
$ws.Cells.Item(2, 20).Value = 0.06198533347787069

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp7"
$ws.Cells.Item(3, 3).Value = "Acvr1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1301303333333333
$ws.Cells.Item(3, 8).Value = 0.390391
$ws.Cells.Item(3, 9).Value = 0.2986126887311651
$ws.Cells.Item(3, 10).Value = 0.2986126887311651
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 13.148327
$ws.Cells.Item(3, 14).Value = 39.444981
$ws.Cells.Item(3, 15).Value = 0.3214784855238645
$ws.Cells.Item(3, 16).Value = 0.3214784855238645
$ws.Cells.Item(3, 17).Value = 1.710996175285667
$ws.Cells.Item(3, 18).Value = 15.398965577571
$ws.Cells.Item(3, 19).Value = 0.09599755493150411
$ws.Cells.Item(3, 20).Value = 0.09599755493150411

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bmp7"
$ws.Cells.Item(4, 3).Value = "Acvr1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1301303333333333
$ws.Cells.Item(4, 8).Value = 0.390391
$ws.Cells.Item(4, 9).Value = 0.2986126887311651
$ws.Cells.Item(4, 10).Value = 0.2986126887311651
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.630791333333334
$ws.Cells.Item(4, 14).Value = 16.892374
$ws.Cells.Item(4, 15).Value = 0.1376736576555254
$ws.Cells.Item(4, 16).Value = 0.1376736576555254
$ws.Cells.Item(4, 17).Value = 0.7327367531371111
$ws.Cells.Item(4, 18).Value = 6.594630778233999
$ws.Cells.Item(4, 19).Value = 0.0411111010799704
$ws.Cells.Item(4, 20).Value = 0.0411111010799704

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Bmp7"
$ws.Cells.Item(5, 3).Value = "Acvr1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1301303333333333
$ws.Cells.Item(5, 8).Value = 0.390391
$ws.Cells.Item(5, 9).Value = 0.2986126887311651
$ws.Cells.Item(5, 10).Value = 0.2986126887311651
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.738585333333333
$ws.Cells.Item(5, 14).Value = 20.215756
$ws.Cells.Item(5, 15).Value = 0.1647593802263456
$ws.Cells.Item(5, 16).Value = 0.1647593802263456
$ws.Cells.Item(5, 17).Value = 0.8768943556217776
$ws.Cells.Item(5, 18).Value = 7.892049200595999
$ws.Cells.Item(5, 19).Value = 0.0491992415230694
$ws.Cells.Item(5, 20).Value = 0.0491992415230694

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Bmp7"
$ws.Cells.Item(6, 3).Value = "Acvr1"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.1301303333333333
$ws.Cells.Item(6, 8).Value = 0.390391
$ws.Cells.Item(6, 9).Value = 0.2986126887311651
$ws.Cells.Item(6, 10).Value = 0.2986126887311651
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9376886666666667
$ws.Cells.Item(6, 14).Value = 2.813066
$ws.Cells.Item(6, 15).Value = 0.02292662271427322
$ws.Cells.Item(6, 16).Value = 0.02292662271427321
$ws.Cells.Item(6, 17).Value = 0.1220217387562222
$ws.Cells.Item(6, 18).Value = 1.098195648806
$ws.Cells.Item(6, 19).Value = 0.006846180452234127
$ws.Cells.Item(6, 20).Value = 0.006846180452234126

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Bmp7"
$ws.Cells.Item(7, 3).Value = "Acvr1"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.1301303333333333
$ws.Cells.Item(7, 8).Value = 0.390391
$ws.Cells.Item(7, 9).Value = 0.2986126887311651
$ws.Cells.Item(7, 10).Value = 0.2986126887311651
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.954327333333333
$ws.Cells.Item(7, 14).Value = 17.862982
$ws.Cells.Item(7, 15).Value = 0.1455841593712531
$ws.Cells.Item(7, 16).Value = 0.1455841593712531
$ws.Cells.Item(7, 17).Value = 0.7748386006624443
$ws.Cells.Item(7, 18).Value = 6.973547405961999
$ws.Cells.Item(7, 19).Value = 0.04347327726651634
$ws.Cells.Item(7, 20).Value = 0.04347327726651634

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Bmp7"
$ws.Cells.Item(8, 3).Value = "Acvr1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.2842186666666667
$ws.Cells.Item(8, 8).Value = 0.852656
$ws.Cells.Item(8, 9).Value = 0.6522022811047395
$ws.Cells.Item(8, 10).Value = 0.6522022811047393
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.489835333333334
$ws.Cells.Item(8, 14).Value = 25.469506
$ws.Cells.Item(8, 15).Value = 0.2075776945087381
$ws.Cells.Item(8, 16).Value = 0.2075776945087381
$ws.Cells.Item(8, 17).Value = 2.412969678659556
$ws.Cells.Item(8, 18).Value = 21.716727107936
$ws.Cells.Item(8, 19).Value = 0.1353826458650618
$ws.Cells.Item(8, 20).Value = 0.1353826458650617

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Bmp7"
$ws.Cells.Item(9, 3).Value = "Acvr1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.2842186666666667
$ws.Cells.Item(9, 8).Value = 0.852656
$ws.Cells.Item(9, 9).Value = 0.6522022811047395
$ws.Cells.Item(9, 10).Value = 0.6522022811047393
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 13.148327
$ws.Cells.Item(9, 14).Value = 39.444981
$ws.Cells.Item(9, 15).Value = 0.3214784855238645
$ws.Cells.Item(9, 16).Value = 0.3214784855238645
$ws.Cells.Item(9, 17).Value = 3.736999968837333
$ws.Cells.Item(9, 18).Value = 33.632999719536
$ws.Cells.Item(9, 19).Value = 0.2096690015847614
$ws.Cells.Item(9, 20).Value = 0.2096690015847613

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Bmp7"
$ws.Cells.Item(10, 3).Value = "Acvr1"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.2842186666666667
$ws.Cells.Item(10, 8).Value = 0.852656
$ws.Cells.Item(10, 9).Value = 0.6522022811047395
$ws.Cells.Item(10, 10).Value = 0.6522022811047393
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.630791333333334
$ws.Cells.Item(10, 14).Value = 16.892374
$ws.Cells.Item(10, 15).Value = 0.1376736576555254
$ws.Cells.Item(10, 16).Value = 0.1376736576555254
$ws.Cells.Item(10, 17).Value = 1.600376005038222
$ws.Cells.Item(10, 18).Value = 14.403384045344
$ws.Cells.Item(10, 19).Value = 0.08979107357096669
$ws.Cells.Item(10, 20).Value = 0.08979107357096663

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Bmp7"
$ws.Cells.Item(11, 3).Value = "Acvr1"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.2842186666666667
$ws.Cells.Item(11, 8).Value = 0.852656
$ws.Cells.Item(11, 9).Value = 0.6522022811047395
$ws.Cells.Item(11, 10).Value = 0.6522022811047393
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 6.738585333333333
$ws.Cells.Item(11, 14).Value = 20.215756
$ws.Cells.Item(11, 15).Value = 0.1647593802263456
$ws.Cells.Item(11, 16).Value = 0.1647593802263456
$ws.Cells.Item(11, 17).Value = 1.915231738659555
$ws.Cells.Item(11, 18).Value = 17.237085647936
$ws.Cells.Item(11, 19).Value = 0.1074564436170257
$ws.Cells.Item(11, 20).Value = 0.1074564436170257

# Row 12
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Bmp7"
$ws.Cells.Item(12, 3).Value = "Acvr1"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.2842186666666667
$ws.Cells.Item(12, 8).Value = 0.852656
$ws.Cells.Item(12, 9).Value = 0.6522022811047395
$ws.Cells.Item(12, 10).Value = 0.6522022811047393
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.9376886666666667
$ws.Cells.Item(12, 14).Value = 2.813066
$ws.Cells.Item(12, 15).Value = 0.02292662271427322
$ws.Cells.Item(12, 16).Value = 0.02292662271427321
$ws.Cells.Item(12, 17).Value = 0.2665086225884445
$ws.Cells.Item(12, 18).Value = 2.398577603296
$ws.Cells.Item(12, 19).Value = 0.01495279563227673
$ws.Cells.Item(12, 20).Value = 0.01495279563227672

# Row 13
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Bmp7"
$ws.Cells.Item(13, 3).Value = "Acvr1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.2842186666666667
$ws.Cells.Item(13, 8).Value = 0.852656
$ws.Cells.Item(13, 9).Value = 0.6522022811047395
$ws.Cells.Item(13, 10).Value = 0.6522022811047393
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 5.954327333333333
$ws.Cells.Item(13, 14).Value = 17.862982
$ws.Cells.Item(13, 15).Value = 0.1455841593712531
$ws.Cells.Item(13, 16).Value = 0.1455841593712531
$ws.Cells.Item(13, 17).Value = 1.692330975576889
$ws.Cells.Item(13, 18).Value = 15.230978780192
$ws.Cells.Item(13, 19).Value = 0.09495032083464723
$ws.Cells.Item(13, 20).Value = 0.0949503208346472

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Bmp7"
$ws.Cells.Item(14, 3).Value = "Acvr1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.021434
$ws.Cells.Item(14, 8).Value = 0.064302
$ws.Cells.Item(14, 9).Value = 0.04918503016409543
$ws.Cells.Item(14, 10).Value = 0.04918503016409543
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 8.489835333333334
$ws.Cells.Item(14, 14).Value = 25.469506
$ws.Cells.Item(14, 15).Value = 0.2075776945087381
$ws.Cells.Item(14, 16).Value = 0.2075776945087381
$ws.Cells.Item(14, 17).Value = 0.1819711305346667
$ws.Cells.Item(14, 18).Value = 1.637740174812
$ws.Cells.Item(14, 19).Value = 0.01020971516580567
$ws.Cells.Item(14, 20).Value = 0.01020971516580567

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Bmp7"
$ws.Cells.Item(15, 3).Value = "Acvr1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.021434
$ws.Cells.Item(15, 8).Value = 0.064302
$ws.Cells.Item(15, 9).Value = 0.04918503016409543
$ws.Cells.Item(15, 10).Value = 0.04918503016409543
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 13.148327
$ws.Cells.Item(15, 14).Value = 39.444981
$ws.Cells.Item(15, 15).Value = 0.3214784855238645
$ws.Cells.Item(15, 16).Value = 0.3214784855238645
$ws.Cells.Item(15, 17).Value = 0.281821240918
$ws.Cells.Item(15, 18).Value = 2.536391168262
$ws.Cells.Item(15, 19).Value = 0.01581192900759899
$ws.Cells.Item(15, 20).Value = 0.01581192900759899

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Bmp7"
$ws.Cells.Item(16, 3).Value = "Acvr1"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.021434
$ws.Cells.Item(16, 8).Value = 0.064302
$ws.Cells.Item(16, 9).Value = 0.04918503016409543
$ws.Cells.Item(16, 10).Value = 0.04918503016409543
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 5.630791333333334
$ws.Cells.Item(16, 14).Value = 16.892374
$ws.Cells.Item(16, 15).Value = 0.1376736576555254
$ws.Cells.Item(16, 16).Value = 0.1376736576555254
$ws.Cells.Item(16, 17).Value = 0.1206903814386667
$ws.Cells.Item(16, 18).Value = 1.086213432948
$ws.Cells.Item(16, 19).Value = 0.006771483004588367
$ws.Cells.Item(16, 20).Value = 0.006771483004588365

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Bmp7"
$ws.Cells.Item(17, 3).Value = "Acvr1"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.021434
$ws.Cells.Item(17, 8).Value = 0.064302
$ws.Cells.Item(17, 9).Value = 0.04918503016409543
$ws.Cells.Item(17, 10).Value = 0.04918503016409543
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 6.738585333333333
$ws.Cells.Item(17, 14).Value = 20.215756
$ws.Cells.Item(17, 15).Value = 0.1647593802263456
$ws.Cells.Item(17, 16).Value = 0.1647593802263456
$ws.Cells.Item(17, 17).Value = 0.1444348380346666
$ws.Cells.Item(17, 18).Value = 1.299913542312
$ws.Cells.Item(17, 19).Value = 0.008103695086250475
$ws.Cells.Item(17, 20).Value = 0.008103695086250475

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Bmp7"
$ws.Cells.Item(18, 3).Value = "Acvr1"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.021434
$ws.Cells.Item(18, 8).Value = 0.064302
$ws.Cells.Item(18, 9).Value = 0.04918503016409543
$ws.Cells.Item(18, 10).Value = 0.04918503016409543
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.9376886666666667
$ws.Cells.Item(18, 14).Value = 2.813066
$ws.Cells.Item(18, 15).Value = 0.02292662271427322
$ws.Cells.Item(18, 16).Value = 0.02292662271427321
$ws.Cells.Item(18, 17).Value = 0.02009841888133333
$ws.Cells.Item(18, 18).Value = 0.180885769932
$ws.Cells.Item(18, 19).Value = 0.001127646629762364
$ws.Cells.Item(18, 20).Value = 0.001127646629762364

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Bmp7"
$ws.Cells.Item(19, 3).Value = "Acvr1"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.021434
$ws.Cells.Item(19, 8).Value = 0.064302
$ws.Cells.Item(19, 9).Value = 0.04918503016409543
$ws.Cells.Item(19, 10).Value = 0.04918503016409543
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 5.954327333333333
$ws.Cells.Item(19, 14).Value = 17.862982
$ws.Cells.Item(19, 15).Value = 0.1455841593712531
$ws.Cells.Item(19, 16).Value = 0.1455841593712531
$ws.Cells.Item(19, 17).Value = 0.1276250520626666
$ws.Cells.Item(19, 18).Value = 1.148625468564
$ws.Cells.Item(19, 19).Value = 0.007160561270089561
$ws.Cells.Item(19, 20).Value = 0.007160561270089561
